$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

Set-TextValue "D2" "332.04"
Set-TextValue "E2" "0.27%"
Set-TextValue "D3" "41.27"
Set-TextValue "E3" "1.33%"
Set-TextValue "D4" "5.713"
Set-TextValue "E4" "-4.84%"
Set-TextValue "D5" "0.08072"
Set-TextValue "E5" "-1.68%"
Set-TextValue "D6" "2.023"
Set-TextValue "E6" "1.97%"
Set-TextValue "D7" "8.747"
Set-TextValue "E7" "-0.65%"
Set-TextValue "D8" "4.535"
Set-TextValue "E8" "-0.84%"
Set-TextValue "D9" "2.949"
Set-TextValue "E9" "0.25%"
Set-TextValue "D10" "0.9229"
Set-TextValue "E10" "-2.90%"
Set-TextValue "D11" "0.1259"
Set-TextValue "E11" "-6.82%"
Set-TextValue "D12" "0.1941"
Set-TextValue "E12" "-2.94%"
Set-TextValue "D13" "8.836"
Set-TextValue "E13" "-11.61%"
Set-TextValue "D14" "0.09343"
Set-TextValue "E14" "0.99%"
Set-TextValue "D15" "0.03688"
Set-TextValue "E15" "5.45%"
Set-TextValue "D16" "0.1050"
Set-TextValue "E16" "9.08%"
Set-TextValue "D17" "0.001302"
Set-TextValue "E17" "-0.72%"
Set-TextValue "D18" "0.006237"
Set-TextValue "E18" "1.16%"
Set-TextValue "D19" "3.362"
Set-TextValue "E19" "0.28%"
Set-TextValue "E20" "-1.62%"
Set-TextValue "D21" "0.1418"
Set-TextValue "E21" "-1.86%"
Set-TextValue "E22" "9.24%"
Set-TextValue "D23" "0.04432"
Set-TextValue "E23" "-0.07%"
Set-TextValue "D24" "0.001262"
Set-TextValue "E24" "0.64%"
Set-TextValue "D25" "0.004310"
Set-TextValue "E25" "-2.98%"
Set-TextValue "D26" "0.0001245"
Set-TextValue "E26" "4.63%"
Set-TextValue "D39" "0.02892"
Set-TextValue "E39" "14.98%"
Set-TextValue "D40" "0.05485"
Set-TextValue "E40" "3.53%"
Set-TextValue "D41" "0.007788"
Set-TextValue "E41" "4.38%"
Set-TextValue "D42" "0.009921"
Set-TextValue "E42" "10.00%"
Set-TextValue "D43" "0.1420"
Set-TextValue "E43" "-2.23%"
Set-TextValue "D44" "0.002240"
Set-TextValue "E44" "9.53%"
Set-TextValue "D45" "0.01114"
Set-TextValue "E45" "5.21%"
Set-TextValue "D46" "0.00006824"
Set-TextValue "E46" "0.54%"
Set-TextValue "E47" "0.74%"
Set-TextValue "D48" "0.002283"
Set-TextValue "E48" "27.27%"
Set-TextValue "D49" "0.003023"
Set-TextValue "E49" "-12.75%"
Set-TextValue "D50" "0.00002108"
Set-TextValue "E50" "0.74%"
Set-TextValue "D51" "0.0002008"
Set-TextValue "E51" "0.74%"
